# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns
# to the latest scraped values, per the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.496.40'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '1.568.79'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").Value = "'" + '208.13'
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("D8").Value = "'" + '22.01'
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("D9").Value = "'" + '0.249'
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("D10").Value = "'" + '0.0590'
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '1.791.72'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Value = '1.573.42'
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").Value = "'" + '3.83'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("E15").Value = '  -1.82%  '
$ws.Range("D16").Value = "'" + '63.37'
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("D17").Value = '27.466.97'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").Value = "'" + '214.19'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").Value = '0.0₃0690'
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("D20").Value = "'" + '7.27'
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("D22").Value = "'" + '4.13'
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("D23").Value = "'" + '9.56'
$ws.Range("E23").Value = '  +1.13%  '
$ws.Range("E24").Value = '  +1.83%  '
$ws.Range("D25").Value = "'" + '153.24'
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("D26").Value = "'" + '6.79'
$ws.Range("E26").Value = '  +1.80%  '
$ws.Range("E27").Value = '  -0.39%  '
$ws.Range("D28").Value = "'" + '15.03'
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("E31").Value = '  +1.78%  '
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("D33").Value = '1.362.64'
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("E34").Value = '  +0.97%  '
$ws.Range("E35").Value = '  +3.00%  '
$ws.Range("D36").Value = "'" + '0.977'
$ws.Range("E36").Value = '  +0.62%  '
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("E38").Value = '  +2.01%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("D40").Value = "'" + '0.822'
$ws.Range("E40").Value = '  +2.37%  '
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  +2.15%  '
$ws.Range("D44").Value = "'" + '64.23'
$ws.Range("E44").Value = '  +1.88%  '
$ws.Range("E45").Value = '  +0.99%  '
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("D47").Value = '1.704.67'
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").Value = "'" + '85.55'
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("D49").Value = '0.0₇0988'
$ws.Range("E49").Value = '  +2.55%  '
$ws.Range("D50").Value = "'" + '0.0955'
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("E51").Value = '  -0.02%  '
